function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "M4" "22"
Set-TextValue $ws "M5" "22"
Set-TextValue $ws "M7" "19"
Set-TextValue $ws "M8" "11"
Set-TextValue $ws "B9" "Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA; University of Ulsan College of Medicine, Asan Medical Center, Seoul, Korea; US Army Institute of Surgical Research, San Antonio, TX, USA; University of Pittsburgh School of Medicine, Pittsburgh, PA, USA; University of Pittsburgh School of Medicine, Pittsburgh, PA, USA; University of Pittsburgh School of Medicine, Pittsburgh, PA, USA; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, USA"
Set-TextValue $ws "M9" "8"
Set-TextValue $ws "M10" "6"
Set-TextValue $ws "A11" "Nancy G. Casanova, Vivian Reyes-Hernon, Taylor Gregory, Belinda Sun, Tadeo Bermudez, Matthew Hufford, Radu C. Oita, Sara M. Camp, Gabriela Hernández-Molina, Jorge Rojas Serrano, Xiaoguang Sun, Jocelyn Fimbres, Mehdi Mirsaeidi, Saad Sammani, Christian Bime, Joe G.N. Garcia"
Set-TextValue $ws "A12" "Jacqueline C. Stocking, Christiana Drake, Janet Aldrich, Michael Ong, Alpesh Amin, Rebecca A. Marmor, Laura N Godat, Maxime Cannesson, Michael A. Gropper, Patrick S. Romano, Christian Sandrock, Christian Bime, Ivo Abraham, Garth H. Utter"
Set-TextValue $ws "B11" "Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Pathology, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Instituto Nacional de Ciencias Médicas y Nutrición Salvador Zubirán, México City, Mexico; Instituto Nacional de Enfermedades Respiratorias, México City, Mexico; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Instituto Nacional de Ciencias Médicas y Nutrición Salvador Zubirán, México City, Mexico; Department of Medicine, College of Medicine, University of Florida, Jacksonville, FL, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States; Department of Medicine, University of Arizona Health Sciences, Tucson, AZ, United States"
Set-TextValue $ws "B12" "Department of Internal Medicine, Division of Pulmonary, Critical Care and Sleep Medicine, University of California Davis, 4150 V Street, Suite 3400, Sacramento, CA, 95817, USA. jcstocking@ucdavis.edu.; Department of Statistics, University of California Davis, Davis, USA; Department of Anesthesia and Perioperative Care, University of California San Francisco, San Francisco, USA; Department of Medicine, University of California Los Angeles, Los Angeles, USA; VA Greater Los Angeles Healthcare System, Los Angeles, USA; Department of Medicine, University of California Irvine, Irvine, USA; Department of Surgery, University of California San Diego, San Diego, USA; Department of Surgery, University of California San Diego, San Diego, USA; Department of Anesthesiology and Perioperative Medicine, University of California Los Angeles, Los Angeles, USA; Department of Anesthesia and Perioperative Care, University of California San Francisco, San Francisco, USA; Center for Healthcare Policy and Research, University of California Davis, Sacramento, USA; Department of Internal Medicine, Division of Pulmonary, Critical Care and Sleep Medicine, University of California Davis, Sacramento, USA; College of Medicine, University of Arizona Health Sciences, Tucson, USA; Center for Health Outcomes and PharmacoEconomic Research, University of Arizona, Tucson, USA; Center for Healthcare Policy and Research, University of California Davis, Sacramento, USA; Department of Surgery, Outcomes Research Group, University of California Davis, Sacramento, USA"
Set-TextValue $ws "C11" "https://openalex.org/W4307371954"
Set-TextValue $ws "C12" "https://openalex.org/W4280488836"
Set-TextValue $ws "D11" "Biochemical and genomic identification of novel biomarkers in progressive sarcoidosis: HBEGF, eNAMPT, and ANG-2"
Set-TextValue $ws "D12" "Outcomes and risk factors for delayed-onset postoperative respiratory failure: a multi-center case-control study by the University of California Critical Care Research Collaborative (UC3RC)"
Set-TextValue $ws "E11" "2022-10-25"
Set-TextValue $ws "E12" "2022-05-14"
Set-TextValue $ws "F11" "Frontiers in Medicine"
Set-TextValue $ws "F12" "BMC Anesthesiology"
Set-TextValue $ws "G11" "Frontiers Media"
Set-TextValue $ws "G12" "BioMed Central"
Set-TextValue $ws "H11" "https://doi.org/10.3389/fmed.2022.1012827"
Set-TextValue $ws "H12" "https://doi.org/10.1186/s12871-022-01681-x"
Set-TextValue $ws "O11" "https://pubmed.ncbi.nlm.nih.gov/36388923"
Set-TextValue $ws "O12" "https://pubmed.ncbi.nlm.nih.gov/35568812"
Set-TextValue $ws "P11" "https://doi.org/10.3389/fmed.2022.1012827"
Set-TextValue $ws "P12" "https://doi.org/10.1186/s12871-022-01681-x"
Set-TextValue $ws "M12" "3"
Set-TextValue $ws "A14" "Christian Bime, Juan C. Celedón"
Set-TextValue $ws "A15" "Jacqueline C. Stocking, Christiana Drake, Janet Aldrich, Michael Ong, Alpesh Amin, Rebecca A. Marmor, Laura N Godat, Maxine Cannesson, Michael A. Gropper, Patrick S. Romano, Christian Sandrock, Christian Bime, Ivo Abraham, Garth H. Utter"
Set-TextValue $ws "B14" "; "
Set-TextValue $ws "B15" "University of California Davis Medical Center; University of California Davis; University of California San Francisco; University of California Los Angeles; University of California Irvine; University of California San Diego; University of California San Diego; University of California Los Angeles; University of California San Francisco; University of California Davis; University of California Davis; The University of Arizona Health Sciences; The University of Arizona; University of California Davis"
Set-TextValue $ws "C14" "https://openalex.org/W4253274986"
Set-TextValue $ws "C15" "https://openalex.org/W4220969699"
Set-TextValue $ws "D14" "Respiratory Health in Migrants and Refugees"
Set-TextValue $ws "D15" "Outcomes and Risk Factors for Delayed-Onset Postoperative Respiratory Failure: A Multi-Center Case-Control Study by the University of California Critical Care Research Collaborative (UC3RC)"
Set-TextValue $ws "E14" "2022-01-01"
Set-TextValue $ws "E15" "2022-03-09"
Set-TextValue $ws "F14" "Elsevier eBooks"
Set-TextValue $ws "F15" "Research Square (Research Square)"
Set-TextValue $ws "G14" "N/A"
Set-TextValue $ws "G15" "Research Square (United States)"
Set-TextValue $ws "H14" "https://doi.org/10.1016/b978-0-08-102723-3.00063-9"
Set-TextValue $ws "H15" "https://doi.org/10.21203/rs.3.rs-824536/v1"
Set-TextValue $ws "I14" "N/A"
Set-TextValue $ws "I15" "cc-by"
Set-TextValue $ws "J14" "N/A"
Set-TextValue $ws "J15" "submittedVersion"
Set-TextValue $ws "K14" "closed"
Set-TextValue $ws "K15" "green"
Set-TextValue $ws "P14" "https://doi.org/10.1016/b978-0-08-102723-3.00063-9"
Set-TextValue $ws "P15" "https://doi.org/10.21203/rs.3.rs-824536/v1"
Set-TextValue $ws "Q14" "book-chapter"
Set-TextValue $ws "Q15" "article"
Set-TextValue $ws "A16" "Mladen Jergović, Makiko Watanabe, Ruchika Bhat, Christopher P Coplen, Sandip Ashok Sonar, Rachel Wong, Yvonne Castaneda, Lisa Davidson, Mrinalini Kala, Rachel C. Wilson, Homer L. Twigg, Kenneth S. Knox, Heidi E Erickson, Craig Weinkauf, Christian Bime, Billie Bixby, Sairam Parthasarathy, Jarrod Mosier, Bonnie LaFleur, Deepta Bhattacharya, Janko Nikolich‐Žugich"
